# Insert a new data row at row 27 (shifts existing rows 27..126 down to 28..127)
# and populate it with the new weekly price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(27).Insert()

$ws.Cells.Item(27, 1).Value  = 3
$ws.Cells.Item(27, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(27, 3).Value  = "Coquimbo"
$ws.Cells.Item(27, 4).Value  = 44701
$ws.Cells.Item(27, 5).Value  = 5
$ws.Cells.Item(27, 6).Value  = "Fruta"
$ws.Cells.Item(27, 7).Value  = 100107
$ws.Cells.Item(27, 8).Value  = "Otros"
$ws.Cells.Item(27, 9).Value  = 100107011
$ws.Cells.Item(27, 10).Value = "Tuna"
$ws.Cells.Item(27, 11).Value = "Sin especificar"
$ws.Cells.Item(27, 12).Value = "Primera"
$ws.Cells.Item(27, 13).Value = 65
$ws.Cells.Item(27, 14).Value = 17000
$ws.Cells.Item(27, 15).Value = 17000
$ws.Cells.Item(27, 16).Value = 17000
$ws.Cells.Item(27, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(27, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(27, 19).Value = 850
$ws.Cells.Item(27, 20).Value = 20
